$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price/Volume columns so that
# numeric-looking strings (e.g. "15.61", "1.659", "8.600") are not
# silently coerced into floating point numbers, and multi-dot
# strings (e.g. "28.199.91") keep their original literal text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.199.91"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "1.878.50"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "316.05"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4314"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("D8").Value = "0.3696"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "0.07421"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").Value = "0.8854"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").Value = "21.16"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").Value = "1.925.85"
$ws.Range("E12").Value = "  +4.69%  "
$ws.Range("D13").Value = "5.474"
$ws.Range("E13").Value = "  +2.37%  "
$ws.Range("D14").Value = "6.627"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "0.06989"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "81.22"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "0.000009139"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "15.61"
$ws.Range("D21").Value = "28.321.77"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").Value = "5.095"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").Value = "10.95"
$ws.Range("E23").Value = "  +2.93%  "
$ws.Range("D24").Value = "2.109.55"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").Value = "1.983"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "154.40"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "5.424"
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("D29").Value = "118.36"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").Value = "1.894"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").Value = "0.08984"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "0.7937"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("D33").Value = "4.718"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").Value = "1.178"
$ws.Range("E34").Value = "  +6.91%  "
$ws.Range("D35").Value = "2.958"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "0.05484"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("D39").Value = "0.01971"
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "2.899"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").Value = "0.5175"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "6.887"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "8.600"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("D45").Value = "10.57"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").Value = "0.06579"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "0.4774"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "105.97"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "1.002"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "1.659"
$ws.Range("D51").Value = "1.856"
$ws.Range("E51").Value = "  +5.79%  "
